# Trade #37 closed at 2026-02-16 22:55:48 - base_strategy UP +0.000%
#
# Appends a new trade-log row (row 38) to the "All Trades" and
# "base_strategy" worksheets, mirroring the existing row layout:
#   A: Trade #   B: Date   C: Time   D: Strategy   E: Side
#   F: Entry Price   G: Exit Price   H: Status   I: P&L %   J: P&L $
#   K: Capital After   L: Entry Slippage (bps)   M: Exit Slippage (bps)
#   N: Confidence   O: Entry Reason   P: Exit Reason   Q: Duration (min)

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 38

    $ws.Cells.Item($row, 1).Value = 37

    # Force the date/time-shaped text into a plain text cell so it is
    # stored verbatim ("2026-02-16") instead of being auto-parsed into a
    # date serial number.
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-16"

    $ws.Cells.Item($row, 3).Value = "22:55:48"
    $ws.Cells.Item($row, 4).Value = "base_strategy"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 49.999998
    # G (Exit Price) stays blank - trade is still OPEN
    $ws.Cells.Item($row, 8).Value = "OPEN"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    # P (Exit Reason) stays blank - trade is still OPEN
    $ws.Cells.Item($row, 17).Value = 0
}
